# Generate Report for Handback
#
# This script applies the "handback" update to the localization-status
# workbook:
#   - Overview sheet: status text for zh-cn / de-de flips from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - zh-cn / de-de detail sheets: status text updated the same way, the
#     "Latest Target File" column gets a hyperlink to the source .md file,
#     the "Latest Handback File" column is filled in with the generated
#     xlf file name, and "Latest Handback DateTime" gets the handback
#     timestamp
#   - column widths on the affected columns are widened to fit the newly
#     populated long file names

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"
$targetFileName = "1f6f1c98-d1ec-4fbc-b0da-acc6453eadbd.md"

# ---- Overview sheet ------------------------------------------------------
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---- zh-cn sheet -----------------------------------------------------------
$zhcn.Range("C2").Value = $statusHandedBack

$zhcn.Range("I2").Value = $targetFileName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55c29ae1093d70995a53da103a66d55973191d06/e2e/1f6f1c98-d1ec-4fbc-b0da-acc6453eadbd.md", [Type]::Missing, $targetFileName, $targetFileName) | Out-Null

$zhcn.Range("J2").Value = "1f6f1c98-d1ec-4fbc-b0da-acc6453eadbd.d6b7a3bab36f9d9e206fe0a250bea01cba82202c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-20 07:02:56"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet -----------------------------------------------------------
$dede.Range("C2").Value = $statusHandedBack

$dede.Range("I2").Value = $targetFileName
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55c29ae1093d70995a53da103a66d55973191d06/e2e/1f6f1c98-d1ec-4fbc-b0da-acc6453eadbd.md", [Type]::Missing, $targetFileName, $targetFileName) | Out-Null

$dede.Range("J2").Value = "1f6f1c98-d1ec-4fbc-b0da-acc6453eadbd.d6b7a3bab36f9d9e206fe0a250bea01cba82202c.de-de.xlf"
$dede.Range("K2").Value = "2016-08-20 07:03:05"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
